$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.63
$ws.Range("I2").Value = 5.8
$ws.Range("J2").Value = 4.6
$ws.Range("L2").Value = 1.32
$ws.Range("N2").Value = 5.2
$ws.Range("O2").Value = 1.22
$ws.Range("P2").Value = 2.42
$ws.Range("Q2").Value = 1.68
$ws.Range("R2").Value = 1.56
$ws.Range("S2").Value = 2.72
$ws.Range("T2").Value = 1.74
$ws.Range("U2").Value = 2.26
$ws.Range("V2").Value = 1.2
$ws.Range("X2").Value = 25
$ws.Range("Y2").Value = 24
$ws.Range("Z2").Value = 55
$ws.Range("AB2").Value = 11
$ws.Range("AC2").Value = 10.5
$ws.Range("AE2").Value = 70
$ws.Range("AF2").Value = 11
$ws.Range("AH2").Value = 19
$ws.Range("AL2").Value = 29
$ws.Range("AN2").Value = 7.6
$ws.Range("F3").Value = 30
$ws.Range("G3").Value = 48
$ws.Range("H3").Value = 1.08
$ws.Range("J3").Value = 17
$ws.Range("K3").Value = 19.5
$ws.Range("L3").Value = 1.08
$ws.Range("P3").Value = 6.6
$ws.Range("Q3").Value = 1.13
$ws.Range("R3").Value = 3.25
$ws.Range("S3").Value = 1.36
$ws.Range("T3").Value = 1.96
$ws.Range("U3").Value = 1.87
$ws.Range("X3").Value = 130
$ws.Range("Z3").Value = 21
$ws.Range("AA3").Value = 14.5
$ws.Range("AC3").Value = 1000
$ws.Range("AD3").Value = 26
$ws.Range("AO3").Value = 1.76
$ws.Range("G4").Value = 2.44
$ws.Range("I4").Value = 4.7
$ws.Range("J4").Value = 2.82
$ws.Range("K4").Value = 3.15
$ws.Range("L4").Value = 1.62
$ws.Range("N4").Value = 2.3
$ws.Range("O4").Value = 1.61
$ws.Range("Q4").Value = 2.86
$ws.Range("T4").Value = 2.26
$ws.Range("U4").Value = 1.64
$ws.Range("V4").Value = 1.28
$ws.Range("W4").Value = 1.69
$ws.Range("X4").Value = 7.8
$ws.Range("AB4").Value = 6.8
$ws.Range("G5").Value = 2.38
$ws.Range("H5").Value = 3.05
$ws.Range("L5").Value = 1.21
$ws.Range("P5").Value = 2.78
$ws.Range("S5").Value = 2.1
$ws.Range("U5").Value = 2.84
$ws.Range("W5").Value = 1.73
$ws.Range("J6").Value = 2.74
$ws.Range("N6").Value = 2.22
$ws.Range("P6").Value = 1.4
$ws.Range("F7").Value = 1.44
$ws.Range("H7").Value = 7.8
$ws.Range("L7").Value = 1.33
$ws.Range("N7").Value = 4.2
$ws.Range("O7").Value = 1.25
$ws.Range("T7").Value = 1.92
$ws.Range("U7").Value = 1.9
$ws.Range("X7").Value = 19.5
$ws.Range("Y7").Value = 1000
$ws.Range("Z7").Value = 75
$ws.Range("AD7").Value = 32
$ws.Range("AE7").Value = 130
$ws.Range("AF7").Value = 9.199999999999999
$ws.Range("AH7").Value = 26
$ws.Range("AI7").Value = 290
$ws.Range("AJ7").Value = 13.5
$ws.Range("AK7").Value = 16
$ws.Range("AL7").Value = 36
$ws.Range("AM7").Value = 150
$ws.Range("AO7").Value = 160
$ws.Range("G8").Value = 2.9
$ws.Range("H8").Value = 2.7
$ws.Range("U8").Value = 2.14
$ws.Range("V8").Value = 1.52
$ws.Range("W8").Value = 1.53
$ws.Range("AO8").Value = 1000
$ws.Range("G9").Value = 1.63
$ws.Range("M9").Value = 1.13
$ws.Range("N9").Value = 2.34
$ws.Range("O9").Value = 1.62
$ws.Range("P9").Value = 1.44
$ws.Range("Q9").Value = 2.86
$ws.Range("R9").Value = 1.15
$ws.Range("T9").Value = 2.82
$ws.Range("U9").Value = 1.46
$ws.Range("AF9").Value = 21
$ws.Range("AG9").Value = 40
$ws.Range("F10").Value = 1.69
$ws.Range("G10").Value = 1.7
$ws.Range("J10").Value = 3.65
$ws.Range("K10").Value = 3.7
$ws.Range("S10").Value = 4.9
$ws.Range("U10").Value = 1.69
$ws.Range("X10").Value = 9.6
$ws.Range("AA10").Value = 240
$ws.Range("AG10").Value = 10.5
$ws.Range("AI10").Value = 160
$ws.Range("AM10").Value = 240
$ws.Range("J11").Value = 13
$ws.Range("K11").Value = 13.5
$ws.Range("P11").Value = 3.1
$ws.Range("Q11").Value = 1.43
$ws.Range("U11").Value = 1.48
$ws.Range("Z11").Value = 600
$ws.Range("AB11").Value = 12
$ws.Range("AH11").Value = 1000
$ws.Range("AN11").Value = 2.88
$ws.Range("G12").Value = 610
$ws.Range("I12").Value = 870
$ws.Range("J12").Value = 1.03
$ws.Range("S12").Value = 1.05
$ws.Range("T12").Value = 1.04
$ws.Range("U12").Value = 1.04
$ws.Range("X12").Value = 990
$ws.Range("Y12").Value = 990
$ws.Range("AB12").Value = 990
$ws.Range("AC12").Value = 990
$ws.Range("AD12").Value = 990
$ws.Range("AG12").Value = 990
$ws.Range("AH12").Value = 990
